$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.00563866666667
$ws.Range("H2").Value = 99.01691600000001
$ws.Range("I2").Value = 0.9169150302490913
$ws.Range("J2").Value = 0.9169150302490912
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8063316666666666
$ws.Range("N2").Value = 2.418995
$ws.Range("O2").Value = 0.1277387112198808
$ws.Range("P2").Value = 0.1277387112198808
$ws.Range("Q2").Value = 26.61349163549111
$ws.Range("R2").Value = 239.52142471942
$ws.Range("S2").Value = 0.1171255442621569
$ws.Range("T2").Value = 0.1171255442621569

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.00563866666667
$ws.Range("H3").Value = 99.01691600000001
$ws.Range("I3").Value = 0.9169150302490913
$ws.Range("J3").Value = 0.9169150302490912
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.578098999999999
$ws.Range("N3").Value = 10.734297
$ws.Range("O3").Value = 0.566840884181833
$ws.Range("P3").Value = 0.5668408841818329
$ws.Range("Q3").Value = 118.0974427075613
$ws.Range("R3").Value = 1062.876984368052
$ws.Range("S3").Value = 0.519744926466007
$ws.Range("T3").Value = 0.5197449264660069

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 33.00563866666667
$ws.Range("H4").Value = 99.01691600000001
$ws.Range("I4").Value = 0.9169150302490913
$ws.Range("J4").Value = 0.9169150302490912
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.215895
$ws.Range("N4").Value = 0.647685
$ws.Range("O4").Value = 0.03420199180918047
$ws.Range("P4").Value = 0.03420199180918047
$ws.Range("Q4").Value = 7.12575235994
$ws.Range("R4").Value = 64.13177123946001
$ws.Range("S4").Value = 0.03136032035429388
$ws.Range("T4").Value = 0.03136032035429388

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 33.00563866666667
$ws.Range("H5").Value = 99.01691600000001
$ws.Range("I5").Value = 0.9169150302490913
$ws.Range("J5").Value = 0.9169150302490912
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.712026
$ws.Range("N5").Value = 5.136078
$ws.Range("O5").Value = 0.2712184127891059
$ws.Range("P5").Value = 0.2712184127891059
$ws.Range("Q5").Value = 56.50651154393867
$ws.Range("R5").Value = 508.5586038954481
$ws.Range("S5").Value = 0.2486842391666335
$ws.Range("T5").Value = 0.2486842391666335

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.07465466666666666
$ws.Range("H6").Value = 0.223964
$ws.Range("I6").Value = 0.002073948231580021
$ws.Range("J6").Value = 0.00207394823158002
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8063316666666666
$ws.Range("N6").Value = 2.418995
$ws.Range("O6").Value = 0.1277387112198808
$ws.Range("P6").Value = 0.1277387112198808
$ws.Range("Q6").Value = 0.06019642179777777
$ws.Range("R6").Value = 0.5417677961799999
$ws.Range("S6").Value = 0.0002649234742387828
$ws.Range("T6").Value = 0.0002649234742387827

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.07465466666666666
$ws.Range("H7").Value = 0.223964
$ws.Range("I7").Value = 0.002073948231580021
$ws.Range("J7").Value = 0.00207394823158002
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.578098999999999
$ws.Range("N7").Value = 10.734297
$ws.Range("O7").Value = 0.566840884181833
$ws.Range("P7").Value = 0.5668408841818329
$ws.Range("Q7").Value = 0.2671217881453333
$ws.Range("R7").Value = 2.404096093308
$ws.Range("S7").Value = 0.001175598649336168
$ws.Range("T7").Value = 0.001175598649336167

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.07465466666666666
$ws.Range("H8").Value = 0.223964
$ws.Range("I8").Value = 0.002073948231580021
$ws.Range("J8").Value = 0.00207394823158002
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.215895
$ws.Range("N8").Value = 0.647685
$ws.Range("O8").Value = 0.03420199180918047
$ws.Range("P8").Value = 0.03420199180918047
$ws.Range("Q8").Value = 0.01611756926
$ws.Range("R8").Value = 0.14505812334
$ws.Range("S8").Value = 0.00007093316042916418
$ws.Range("T8").Value = 0.00007093316042916417

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.07465466666666666
$ws.Range("H9").Value = 0.223964
$ws.Range("I9").Value = 0.002073948231580021
$ws.Range("J9").Value = 0.00207394823158002
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.712026
$ws.Range("N9").Value = 5.136078
$ws.Range("O9").Value = 0.2712184127891059
$ws.Range("P9").Value = 0.2712184127891059
$ws.Range("Q9").Value = 0.1278107303546667
$ws.Range("R9").Value = 1.150296573192
$ws.Range("S9").Value = 0.0005624929475759061
$ws.Range("T9").Value = 0.000562492947575906

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.576299333333333
$ws.Range("H10").Value = 7.728898
$ws.Range("I10").Value = 0.07157103078692272
$ws.Range("J10").Value = 0.0715710307869227
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8063316666666666
$ws.Range("N10").Value = 2.418995
$ws.Range("O10").Value = 0.1277387112198808
$ws.Range("P10").Value = 0.1277387112198808
$ws.Range("Q10").Value = 2.077351735278889
$ws.Range("R10").Value = 18.69616561751
$ws.Range("S10").Value = 0.009142391233399919
$ws.Range("T10").Value = 0.009142391233399917

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.576299333333333
$ws.Range("H11").Value = 7.728898
$ws.Range("I11").Value = 0.07157103078692272
$ws.Range("J11").Value = 0.0715710307869227
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.578098999999999
$ws.Range("N11").Value = 10.734297
$ws.Range("O11").Value = 0.566840884181833
$ws.Range("P11").Value = 0.5668408841818329
$ws.Range("Q11").Value = 9.218254068300666
$ws.Range("R11").Value = 82.96428661470598
$ws.Range("S11").Value = 0.04056938637306446
$ws.Range("T11").Value = 0.04056938637306445

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.576299333333333
$ws.Range("H12").Value = 7.728898
$ws.Range("I12").Value = 0.07157103078692272
$ws.Range("J12").Value = 0.0715710307869227
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.215895
$ws.Range("N12").Value = 0.647685
$ws.Range("O12").Value = 0.03420199180918047
$ws.Range("P12").Value = 0.03420199180918047
$ws.Range("Q12").Value = 0.55621014457
$ws.Range("R12").Value = 5.005891301129999
$ws.Range("S12").Value = 0.002447871808748934
$ws.Range("T12").Value = 0.002447871808748933

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.576299333333333
$ws.Range("H13").Value = 7.728898
$ws.Range("I13").Value = 0.07157103078692272
$ws.Range("J13").Value = 0.0715710307869227
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.712026
$ws.Range("N13").Value = 5.136078
$ws.Range("O13").Value = 0.2712184127891059
$ws.Range("P13").Value = 0.2712184127891059
$ws.Range("Q13").Value = 4.410691442449334
$ws.Range("R13").Value = 39.69622298204401
$ws.Range("S13").Value = 0.01941138137170941
$ws.Range("T13").Value = 0.0194113813717094

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.3398056666666667
$ws.Range("H14").Value = 1.019417
$ws.Range("I14").Value = 0.009439990732406145
$ws.Range("J14").Value = 0.009439990732406145
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.8063316666666666
$ws.Range("N14").Value = 2.418995
$ws.Range("O14").Value = 0.1277387112198808
$ws.Range("P14").Value = 0.1277387112198808
$ws.Range("Q14").Value = 0.2739960695461111
$ws.Range("R14").Value = 2.465964625915
$ws.Range("S14").Value = 0.00120585225008518
$ws.Range("T14").Value = 0.00120585225008518

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.3398056666666667
$ws.Range("H15").Value = 1.019417
$ws.Range("I15").Value = 0.009439990732406145
$ws.Range("J15").Value = 0.009439990732406145
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.578098999999999
$ws.Range("N15").Value = 10.734297
$ws.Range("O15").Value = 0.566840884181833
$ws.Range("P15").Value = 0.5668408841818329
$ws.Range("Q15").Value = 1.215858316094333
$ws.Range("R15").Value = 10.942724844849
$ws.Range("S15").Value = 0.005350972693425408
$ws.Range("T15").Value = 0.005350972693425407

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.3398056666666667
$ws.Range("H16").Value = 1.019417
$ws.Range("I16").Value = 0.009439990732406145
$ws.Range("J16").Value = 0.009439990732406145
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.215895
$ws.Range("N16").Value = 0.647685
$ws.Range("O16").Value = 0.03420199180918047
$ws.Range("P16").Value = 0.03420199180918047
$ws.Range("Q16").Value = 0.07336234440499999
$ws.Range("R16").Value = 0.660261099645
$ws.Range("S16").Value = 0.0003228664857084945
$ws.Range("T16").Value = 0.0003228664857084945

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.3398056666666667
$ws.Range("H17").Value = 1.019417
$ws.Range("I17").Value = 0.009439990732406145
$ws.Range("J17").Value = 0.009439990732406145
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.712026
$ws.Range("N17").Value = 5.136078
$ws.Range("O17").Value = 0.2712184127891059
$ws.Range("P17").Value = 0.2712184127891059
$ws.Range("Q17").Value = 0.5817561362806667
$ws.Range("R17").Value = 5.235805226526001
$ws.Range("S17").Value = 0.002560299303187063
$ws.Range("T17").Value = 0.002560299303187063
